$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values per corrected results

$ws.Range("C2").Value = -0.03622362435067662
$ws.Range("D2").Value = 0.9714307255777141

$ws.Range("C3").Value = 1.089223708455519
$ws.Range("D3").Value = 0.287843262960958

$ws.Range("C4").Value = 0.6965043292860141
$ws.Range("D4").Value = 0.4934071022591615

$ws.Range("C5").Value = 0.08652872502303947
$ws.Range("D5").Value = 0.9318288416960669

$ws.Range("C6").Value = 0.8463249792047624
$ws.Range("D6").Value = 0.4064834492831513

$ws.Range("C7").Value = 0.6076202549619911
$ws.Range("D7").Value = 0.5496606345910346

$ws.Range("C8").Value = 0.09303577641890622
$ws.Range("D8").Value = 0.9267171998618307

$ws.Range("C9").Value = -0.1283304979230291
$ws.Range("D9").Value = 0.8990533407996708

$ws.Range("C10").Value = -0.8272235079504265
$ws.Range("D10").Value = 0.4169932213994252

$ws.Range("C11").Value = -0.580999509617843
$ws.Range("D11").Value = 0.5671458951730961
